# Update "想去人数" (number of people interested) figures that changed
# between scrapes, on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3302
$ws1.Range("F3").Value = 13
$ws1.Range("F5").Value = 1304
$ws1.Range("F6").Value = 314

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3302
$ws4.Range("F3").Value = 13
$ws4.Range("F5").Value = 1304
$ws4.Range("F7").Value = 314
